# Update TPM-derived values in the NATMI LR-pairs sheet.
# Only columns M (Receptor average expression value) and N (Receptor total
# expression value) carry new "raw" TPM-based input numbers; columns O, P,
# Q, R, S, T are recomputed downstream of M/N (and the unchanged ligand
# columns G/H), so every affected cell is written explicitly below with the
# exact values that the new TPM pipeline produced.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("M2").Value = 0.2054816666666667
$ws.Range("N2").Value = 0.616445
$ws.Range("O2").Value = 0.00496424614546655
$ws.Range("P2").Value = 0.004964246145466549
$ws.Range("Q2").Value = 0.002903318962222222
$ws.Range("R2").Value = 0.02612987066
$ws.Range("S2").Value = 0.00496424614546655
$ws.Range("T2").Value = 0.004964246145466549

# Row 3 (only the derived-specificity columns change; M3/N3/Q3/R3 stay put)
$ws.Range("O3").Value = 0.9529850468799925
$ws.Range("P3").Value = 0.9529850468799924
$ws.Range("S3").Value = 0.9529850468799925
$ws.Range("T3").Value = 0.9529850468799924

# Row 4
$ws.Range("M4").Value = 1.712817
$ws.Range("N4").Value = 5.138451
$ws.Range("O4").Value = 0.04138006727350978
$ws.Range("P4").Value = 0.04138006727350978
$ws.Range("Q4").Value = 0.024200962332
$ws.Range("R4").Value = 0.217808660988
$ws.Range("S4").Value = 0.04138006727350978
$ws.Range("T4").Value = 0.04138006727350978

# Row 5
$ws.Range("M5").Value = 0.02775933333333333
$ws.Range("N5").Value = 0.083278
$ws.Range("O5").Value = 0.000670639701031176
$ws.Range("P5").Value = 0.000670639701031176
$ws.Range("Q5").Value = 0.0003922208737777778
$ws.Range("R5").Value = 0.003529987864000001
$ws.Range("S5").Value = 0.000670639701031176
$ws.Range("T5").Value = 0.000670639701031176
